$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (Excel COM uses BGR-packed decimal: R + G*256 + B*65536)
$RED    = 255        # FF0000 -> "Falta esto" style
$GREEN  = 5373729     # 21FF51 -> "Listo" style
$PURPLE = 10498160    # 7030A0 -> new style used on B24

# Target state for column B (status) cells, rows 7-29 (row 19 unchanged)
$rows = @(7,8,9,10,11,12,13,14,15,16,17,18,20,21,22,23,24,25,26,27,28,29)
$values = @{
    7  = "Falta esto"
    8  = "Falta esto"
    9  = "Falta esto"
    10 = "Listo"
    11 = "Listo"
    12 = "Listo"
    13 = "Falta esto"
    14 = "Falta esto"
    15 = "Listo"
    16 = "Cargador de (mas) mapas desde Tiled"
    17 = "Listo"
    18 = "(Falta volver a coger despues de muerto)"
    20 = "Falta esto"
    21 = "Listo (usar otro entregable)"
    22 = "(Subir entregable con los 3 rayos)"
    23 = "Listo"
    24 = "Esto no se ni que es"
    25 = "Ajustarlo"
    26 = "Listo"
    27 = "Listo"
    28 = "Listo"
    29 = "Listo"
}
$colors = @{
    7  = $RED
    8  = $RED
    9  = $RED
    10 = $GREEN
    11 = $GREEN
    12 = $GREEN
    13 = $RED
    14 = $RED
    15 = $GREEN
    16 = $RED
    17 = $GREEN
    18 = $RED
    20 = $RED
    21 = $GREEN
    22 = $RED
    23 = $GREEN
    24 = $PURPLE
    25 = $RED
    26 = $GREEN
    27 = $GREEN
    28 = $GREEN
    29 = $GREEN
}

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = $values[$r]
    $cell.Font.Color = $colors[$r]
}

# Keep the final selection where the author left it
$ws.Range("B18").Select()
